# Update "Förändrad" (changed) date in column C for every data row (2..224)
# from 2023-09-23 (45192) to 2023-10-03 (45202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 2; $i -le 224; $i++) {
    $ws.Cells.Item($i, 3).Value = 45202
}

# Row 2 specific updates: species counts increased because two new species
# (Vågticka, Svartvit taggsvamp) were added to the sighting list.
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 3
$ws.Range("O2").Value = 10
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 16

$ws.Range("R2").Value = "Grangråticka`r`nSpricktaggsvamp`r`nVågticka`r`nDofttaggsvamp`r`nGrantaggsvamp`r`nGul taggsvamp`r`nMotaggsvamp`r`nOrange taggsvamp`r`nSvartvit taggsvamp`r`nTallriska`r`nBlåmossa`r`nBrandticka`r`nDiskvaxskivling`r`nFjällig taggsvamp s.str.`r`nSkarp dropptaggsvamp`r`nTjockfotad fingersvamp"
